$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# Row 102 <-> 103 swap (inverted-odds correction)
Set-Cell $ws 102 2 6007882   # B102
Set-Cell $ws 102 6 "Cimarrones de Sonora FC"   # F102
Set-Cell $ws 102 7 "Club Celaya"   # G102
Set-Cell $ws 102 11 2.1   # K102
Set-Cell $ws 102 12 3.2   # L102
Set-Cell $ws 102 13 3.25   # M102
Set-Cell $ws 102 14 2.2   # N102
Set-Cell $ws 102 15 3   # O102
Set-Cell $ws 102 16 3.6   # P102
Set-Cell $ws 102 17 -0.25   # Q102
Set-Cell $ws 102 18 1.875   # R102
Set-Cell $ws 102 19 1.925   # S102
Set-Cell $ws 102 20 2   # T102
Set-Cell $ws 102 21 1.8   # U102
Set-Cell $ws 102 22 2   # V102
Set-Cell $ws 102 24 2   # X102
Set-Cell $ws 102 26 -0.5   # Z102
Set-Cell $ws 102 27 0.4625   # AA102
Set-Cell $ws 102 28 0   # AB102
Set-Cell $ws 102 29 0   # AC102

# Row 103
Set-Cell $ws 103 2 6007883   # B103
Set-Cell $ws 103 6 "Club Atletico La Paz"   # F103
Set-Cell $ws 103 7 "Dorados"   # G103
Set-Cell $ws 103 11 1.65   # K103
Set-Cell $ws 103 12 3.75   # L103
Set-Cell $ws 103 13 4.5   # M103
Set-Cell $ws 103 14 1.533   # N103
Set-Cell $ws 103 15 4.333   # O103
Set-Cell $ws 103 16 6   # P103
Set-Cell $ws 103 17 -1   # Q103
Set-Cell $ws 103 18 1.825   # R103
Set-Cell $ws 103 19 1.975   # S103
Set-Cell $ws 103 20 2.75   # T103
Set-Cell $ws 103 21 1.95   # U103
Set-Cell $ws 103 22 1.85   # V103
Set-Cell $ws 103 24 3.333   # X103
Set-Cell $ws 103 26 -1   # Z103
Set-Cell $ws 103 27 0.9750000000000001   # AA103
Set-Cell $ws 103 28 -1   # AB103
Set-Cell $ws 103 29 0.8500000000000001   # AC103

# Row 262 <-> 263 swap
Set-Cell $ws 262 2 6924569   # B262
Set-Cell $ws 262 6 "Venados FC"   # F262
Set-Cell $ws 262 7 "Dorados"   # G262
Set-Cell $ws 262 8 4   # H262
Set-Cell $ws 262 10 "H"   # J262
Set-Cell $ws 262 11 1.615   # K262
Set-Cell $ws 262 12 4   # L262
Set-Cell $ws 262 13 4.5   # M262
Set-Cell $ws 262 14 1.5   # N262
Set-Cell $ws 262 15 4.75   # O262
Set-Cell $ws 262 16 5.75   # P262
Set-Cell $ws 262 17 -1.25   # Q262
Set-Cell $ws 262 18 1.925   # R262
Set-Cell $ws 262 19 1.875   # S262
Set-Cell $ws 262 20 3   # T262
Set-Cell $ws 262 21 1.75   # U262
Set-Cell $ws 262 22 1.95   # V262
Set-Cell $ws 262 23 0.5   # W262
Set-Cell $ws 262 25 -1   # Y262
Set-Cell $ws 262 26 0.925   # Z262
Set-Cell $ws 262 27 -1   # AA262
Set-Cell $ws 262 28 0.75   # AB262
Set-Cell $ws 262 29 -1   # AC262

# Row 263
Set-Cell $ws 263 2 6924568   # B263
Set-Cell $ws 263 6 "Atletico Morelia"   # F263
Set-Cell $ws 263 7 "Atlante"   # G263
Set-Cell $ws 263 8 0   # H263
Set-Cell $ws 263 10 "A"   # J263
Set-Cell $ws 263 11 2.4   # K263
Set-Cell $ws 263 12 3   # L263
Set-Cell $ws 263 13 2.875   # M263
Set-Cell $ws 263 14 2.7   # N263
Set-Cell $ws 263 15 3.1   # O263
Set-Cell $ws 263 16 2.8   # P263
Set-Cell $ws 263 17 0   # Q263
Set-Cell $ws 263 18 1.85   # R263
Set-Cell $ws 263 19 1.95   # S263
Set-Cell $ws 263 20 2.25   # T263
Set-Cell $ws 263 21 1.975   # U263
Set-Cell $ws 263 22 1.725   # V263
Set-Cell $ws 263 23 -1   # W263
Set-Cell $ws 263 25 1.8   # Y263
Set-Cell $ws 263 26 -1   # Z263
Set-Cell $ws 263 27 0.95   # AA263
Set-Cell $ws 263 28 -1   # AB263
Set-Cell $ws 263 29 0.7250000000000001   # AC263

# Row 337 odds refresh
Set-Cell $ws 337 14 1.85   # N337
Set-Cell $ws 337 15 3.6   # O337
Set-Cell $ws 337 16 4.333   # P337
Set-Cell $ws 337 18 1.825   # R337
Set-Cell $ws 337 19 1.975   # S337
Set-Cell $ws 337 21 1.825   # U337
Set-Cell $ws 337 22 1.975   # V337

# Row 338 odds refresh
Set-Cell $ws 338 14 2.5   # N338
Set-Cell $ws 338 16 2.9   # P338
Set-Cell $ws 338 17 -0.25   # Q338
Set-Cell $ws 338 18 2.05   # R338
Set-Cell $ws 338 19 1.75   # S338

# Row 339 odds refresh
Set-Cell $ws 339 14 3.5   # N339
Set-Cell $ws 339 15 3.5   # O339
Set-Cell $ws 339 18 2   # R339
Set-Cell $ws 339 19 1.8   # S339
Set-Cell $ws 339 21 2   # U339
Set-Cell $ws 339 22 1.8   # V339

# Row 340 odds refresh
Set-Cell $ws 340 15 3.5   # O340
Set-Cell $ws 340 16 3.8   # P340
Set-Cell $ws 340 18 1.95   # R340
Set-Cell $ws 340 19 1.85   # S340
Set-Cell $ws 340 21 1.85   # U340
Set-Cell $ws 340 22 1.95   # V340

# Row 341 odds refresh
Set-Cell $ws 341 16 4.75   # P341
Set-Cell $ws 341 18 2   # R341
Set-Cell $ws 341 19 1.8   # S341
Set-Cell $ws 341 21 1.875   # U341
Set-Cell $ws 341 22 1.925   # V341

# Row 342 odds refresh
Set-Cell $ws 342 14 1.65   # N342
Set-Cell $ws 342 15 4.2   # O342
Set-Cell $ws 342 16 5   # P342
Set-Cell $ws 342 18 2   # R342
Set-Cell $ws 342 19 1.8   # S342

# New row 343 (new fixture)
Set-Cell $ws 343 1 341   # A343
Set-Cell $ws 343 2 7641680   # B343
Set-Cell $ws 343 3 "Mexico Liga de Expansion"   # C343
Set-Cell $ws 343 4 "Mexico Liga de Expansion"   # D343
Set-Cell $ws 343 5 45347.83680555555   # E343
Set-Cell $ws 343 6 "Club Celaya"   # F343
Set-Cell $ws 343 7 "Oaxaca"   # G343
Set-Cell $ws 343 11 1.333   # K343
Set-Cell $ws 343 12 4.75   # L343
Set-Cell $ws 343 13 7.5   # M343
Set-Cell $ws 343 14 1.363   # N343
Set-Cell $ws 343 15 4.75   # O343
Set-Cell $ws 343 16 7   # P343
Set-Cell $ws 343 17 -1.5   # Q343
Set-Cell $ws 343 18 1.975   # R343
Set-Cell $ws 343 19 1.825   # S343
Set-Cell $ws 343 20 3   # T343
Set-Cell $ws 343 21 1.975   # U343
Set-Cell $ws 343 22 1.825   # V343
Set-Cell $ws 343 23 0   # W343
Set-Cell $ws 343 24 0   # X343
Set-Cell $ws 343 25 0   # Y343
Set-Cell $ws 343 26 0   # Z343
Set-Cell $ws 343 27 0   # AA343

# Copy formatting for the new row's styled cells (A343 bold/border/center, E343 date fmt)
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(343, 1).PasteSpecial(-4122)
$ws.Cells.Item(343, 5).NumberFormat = $ws.Cells.Item(2, 5).NumberFormat
$ws.Application.CutCopyMode = $false

